$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 2924.875
$ws.Range("I32").Value = 3349.75
$ws.Range("K32").Value = 3349.75
$ws.Range("M32").Value = -3023.75

$ws.Range("H69").Value = 6930.6665
$ws.Range("I69").Value = 4013
$ws.Range("K69").Value = 12039
$ws.Range("M69").Value = -11165

$ws.Range("H72").Value = 6930.6665
$ws.Range("I72").Value = 4013
$ws.Range("K72").Value = 36117
$ws.Range("M72").Value = -31749

$ws.Range("H76").Value = 3210.074
$ws.Range("I76").Value = 3115.9421
$ws.Range("K76").Value = 3115.9421
$ws.Range("M76").Value = -2800.9421

$ws.Range("H79").Value = 3210.074
$ws.Range("I79").Value = 3115.9421
$ws.Range("K79").Value = 3115.9421
$ws.Range("M79").Value = -2023.9421

$ws.Range("H112").Value = 2545.0334
$ws.Range("J112").Value = 2650.476
$ws.Range("L112").Value = 7951.428
$ws.Range("N112").Value = -10167.428

$ws.Range("H138").Value = 2336.6316
$ws.Range("J138").Value = 2618.9565
$ws.Range("L138").Value = 7856.869499999999
$ws.Range("N138").Value = -18136.8695

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 147224.64
$ws.Range("I32").Value = 158968.5
$ws.Range("K32").Value = 158968.5
$ws.Range("M32").Value = -158681.5

$ws.Range("H63").Value = 3521.2856
$ws.Range("J63").Value = 6000
$ws.Range("L63").Value = 6000
$ws.Range("N63").Value = -7372

$ws.Range("H66").Value = 3521.2856
$ws.Range("J66").Value = 6000
$ws.Range("L66").Value = 30000
$ws.Range("N66").Value = -36864

$ws.Range("H74").Value = 14696.823
$ws.Range("I74").Value = 1761.7778
$ws.Range("J74").Value = 29248.75
$ws.Range("K74").Value = 1761.7778
$ws.Range("L74").Value = 29248.75
$ws.Range("M74").Value = -887.7778000000001
$ws.Range("N74").Value = -30996.75

$ws.Range("H77").Value = 14696.823
$ws.Range("I77").Value = 1761.7778
$ws.Range("J77").Value = 29248.75
$ws.Range("K77").Value = 8808.889000000001
$ws.Range("L77").Value = 146243.75
$ws.Range("M77").Value = -4440.889000000001
$ws.Range("N77").Value = -154979.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H14").Value = 304
$ws.Range("I14").Value = 304
$ws.Range("K14").Value = 304
$ws.Range("M14").Value = -132

$ws.Range("H17").Value = 454.5
$ws.Range("J17").Value = 454.5
$ws.Range("L17").Value = 454.5
$ws.Range("N17").Value = -798.5

$ws.Range("H99").Value = 8889.538
$ws.Range("I99").Value = 17329
$ws.Range("J99").Value = 1655.7142
$ws.Range("K99").Value = 17329
$ws.Range("L99").Value = 1655.7142
$ws.Range("M99").Value = -15831
$ws.Range("N99").Value = -4651.7142

$ws.Range("H134").Value = 8605.357
$ws.Range("I134").Value = 5220.3706
$ws.Range("K134").Value = 15661.1118
$ws.Range("M134").Value = -13126.1118

$ws.Range("H135").Value = 86333.336
$ws.Range("J135").Value = 86333.336
$ws.Range("L135").Value = 86333.336
$ws.Range("N135").Value = -96473.336

$ws.Range("H137").Value = 67223
$ws.Range("J137").Value = 67223
$ws.Range("L137").Value = 67223
$ws.Range("N137").Value = -77423

$ws.Range("H138").Value = 79997.75
$ws.Range("J138").Value = 79997.75
$ws.Range("L138").Value = 79997.75
$ws.Range("N138").Value = -90277.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3759.4856
$ws.Range("I31").Value = 4375.15
$ws.Range("J31").Value = 2938.6
$ws.Range("K31").Value = 4375.15
$ws.Range("L31").Value = 2938.6
$ws.Range("M31").Value = -4080.15
$ws.Range("N31").Value = -3528.6

$ws.Range("H34").Value = 3759.4856
$ws.Range("I34").Value = 4375.15
$ws.Range("J34").Value = 2938.6
$ws.Range("K34").Value = 4375.15
$ws.Range("L34").Value = 2938.6
$ws.Range("M34").Value = -4173.15
$ws.Range("N34").Value = -3342.6

$ws.Range("H58").Value = 16198.1
$ws.Range("J58").Value = 30197.25
$ws.Range("L58").Value = 30197.25
$ws.Range("N58").Value = -30603.25

$ws.Range("H62").Value = 4084.4443
$ws.Range("I62").Value = 3350
$ws.Range("K62").Value = 3350
$ws.Range("M62").Value = -2726

$ws.Range("H65").Value = 4084.4443
$ws.Range("I65").Value = 3350
$ws.Range("K65").Value = 16750
$ws.Range("M65").Value = -13630

$ws.Range("H86").Value = 19271
$ws.Range("I86").Value = 36633
$ws.Range("J86").Value = 6249.5
$ws.Range("K86").Value = 36633
$ws.Range("L86").Value = 6249.5
$ws.Range("M86").Value = -35510
$ws.Range("N86").Value = -8495.5

$ws.Range("H89").Value = 19271
$ws.Range("I89").Value = 36633
$ws.Range("J89").Value = 6249.5
$ws.Range("K89").Value = 183165
$ws.Range("L89").Value = 31247.5
$ws.Range("M89").Value = -177549
$ws.Range("N89").Value = -42479.5

$ws.Range("H136").Value = 16198.1
$ws.Range("J136").Value = 30197.25
$ws.Range("L136").Value = 90591.75
$ws.Range("N136").Value = -95691.75

$ws.Range("H141").Value = 154860.38
$ws.Range("J141").Value = 162851.06
$ws.Range("L141").Value = 162851.06
$ws.Range("N141").Value = -173211.06

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H98").Value = 2320.2
$ws.Range("I98").Value = 2320.2
$ws.Range("K98").Value = 6960.599999999999
$ws.Range("M98").Value = -5462.599999999999

$ws.Range("H112").Value = 7197.1113
$ws.Range("I112").Value = 1050.8
$ws.Range("J112").Value = 14880
$ws.Range("K112").Value = 3152.4
$ws.Range("L112").Value = 44640
$ws.Range("M112").Value = -2044.4
$ws.Range("N112").Value = -46856

$ws.Range("H113").Value = 1381.909
$ws.Range("I113").Value = 633.5
$ws.Range("J113").Value = 1809.5714
$ws.Range("K113").Value = 1900.5
$ws.Range("L113").Value = 5428.7142
$ws.Range("M113").Value = 269.5
$ws.Range("N113").Value = -9768.7142

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H7").Value = 501916.5
$ws.Range("J7").Value = 2833
$ws.Range("L7").Value = 2833
$ws.Range("N7").Value = -3057

$ws.Range("H8").Value = 501916.5
$ws.Range("J8").Value = 2833
$ws.Range("L8").Value = 2833
$ws.Range("N8").Value = -3111

$ws.Range("H18").Value = 4580.0835
$ws.Range("I18").Value = 3329
$ws.Range("K18").Value = 3329
$ws.Range("M18").Value = -3036

$ws.Range("H132").Value = 7431.8335
$ws.Range("I132").Value = 5679.0605
$ws.Range("J132").Value = 13858.667
$ws.Range("K132").Value = 17037.1815
$ws.Range("L132").Value = 41576.001
$ws.Range("M132").Value = -14507.1815
$ws.Range("N132").Value = -46636.001

$ws.Range("H135").Value = 97915.414
$ws.Range("J135").Value = 99998.74000000001
$ws.Range("L135").Value = 99998.74000000001
$ws.Range("N135").Value = -110138.74

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H23").Value = 7450
$ws.Range("I23").Value = 7450
$ws.Range("K23").Value = 7450
$ws.Range("M23").Value = -7220

$ws.Range("H93").Value = 2698.4285
$ws.Range("I93").Value = 2298.1052
$ws.Range("K93").Value = 2298.1052
$ws.Range("M93").Value = -1050.1052

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H105").Value = 0
$ws.Range("J105").Value = 0
$ws.Range("L105").Value = 0
$ws.Range("N105").ClearContents()
